# Language workbook update:
#  - some level tweaks, added level 4a, added delete all button
#
# This adds four new Key/Value (translation) entries to the "en" sheet:
#   1. level_0_end_2_a   - new line inserted right after level_0_end_2
#   2. level_4_intro_0_4 - new line inserted right after level_4_intro_0_3
#   3. delete_all_title / delete_all_desc / yes / no - new rows appended
#      at the end of the table for the new "Delete All Cubes" confirmation
#      dialog.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert "level_0_end_2_a" row right after "level_0_end_2" (row 67) ---
# This pushes the existing row 68 ("level_0_end_3") and everything below it
# down by one.
$ws.Rows("68:68").Insert()
$ws.Range("A68").Value = "level_0_end_2_a"
$ws.Range("B68").Value = "To simplify this, take the measurement of the length, width, and height; then multiply them. In this case: 4 times 4 times 1 equals 16."

# --- 2. Insert "level_4_intro_0_4" row right after "level_4_intro_0_3" ---
# After the previous insertion, "level_4_intro_0_3" now lives on row 81, so
# the new row goes in at row 82, pushing everything below it down by one.
$ws.Rows("82:82").Insert()
$ws.Range("A82").Value = "level_4_intro_0_4"
$ws.Range("B82").Value = "Simply apply what you've done from the previous level to fill the grid, but with only one type of cube."

# --- 3. Append the new "delete all cubes" confirmation strings ---
# The table now ends at row 87 (end_detail_2), so the new rows go at the
# bottom, 88-91.
$ws.Range("A88").Value = "delete_all_title"
$ws.Range("B88").Value = "Delete All Cubes"

$ws.Range("A89").Value = "delete_all_desc"
$ws.Range("B89").Value = "Are you sure you want to delete all the cubes?"

$ws.Range("A90").Value = "yes"
$ws.Range("B90").Value = "YES"

$ws.Range("A91").Value = "no"
$ws.Range("B91").Value = "NO"

# Match the final selection left behind in the saved workbook.
$ws.Range("B82").Select()
